$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns stay text-formatted so numeric-looking values are not coerced to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.355.74"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.430.78"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "565.37"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").Value = "144.44"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("D9").Value = "2.426.17"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").Value = "26.62"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  -3.35%  "

$ws.Range("D16").Value = "2.873.25"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "62.350.63"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").Value = "2.424.25"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  +2.39%  "

$ws.Range("D21").Value = "324.79"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  +10.10%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "65.04"
$ws.Range("E25").Value = "  -3.18%  "

$ws.Range("D26").Value = "610.81"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("D27").Value = "8.68"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "0.0₃0989"
$ws.Range("E28").Value = "  -3.30%  "

$ws.Range("D29").Value = "2.558.94"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  -5.05%  "

$ws.Range("D33").Value = "1.87"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("D34").Value = "0.136"
$ws.Range("E34").Value = "  -4.23%  "

$ws.Range("D35").Value = "5.06"
$ws.Range("E35").Value = "  +2.14%  "

$ws.Range("E36").Value = "  -2.30%  "

$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("D39").Value = "18.60"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").Value = "5.26"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("D41").Value = "144.97"
$ws.Range("E41").Value = "  -2.22%  "

$ws.Range("D42").Value = "1.76"
$ws.Range("E42").Value = "  -4.27%  "

$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  -1.23%  "

$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").Value = "41.84"
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("D46").Value = "147.17"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").Value = "3.72"
$ws.Range("E47").Value = "  +0.67%  "

$ws.Range("D48").Value = "20.63"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").Value = "0.0526"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("D50").Value = "0.592"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").Value = "0.0230"
$ws.Range("E51").Value = "  -1.43%  "
